# milp + reformulate raw material buying
#
# - Updates computed Amount values on the existing "Production" and
#   "Inventory" sheets (re-solved MILP results).
# - Adds three new sheets: RawMaterial, RawMaterialInventory, PMRunning.

function Set-Num($ws, $r, $c, $v) {
    $ws.Cells.Item($r, $c).Value = [double]$v
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Production sheet - update column E (Amount)
# ---------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Production")

$prodChanges = @(
    @(2, "510.0000000000006"),
    @(3, "0.0"),
    @(20, "170.00000000000068"),
    @(22, "10.0"),
    @(24, "10.0"),
    @(26, "489.9999999999994"),
    @(42, "1.0231815394945443e-12"),
    @(44, "251.9999999999991"),
    @(46, "989.9999999999998"),
    @(48, "990.0000000000002")
)

foreach ($row in $prodChanges) {
    Set-Num $wsProd $row[0] 5 $row[1]
}

# ---------------------------------------------------------------------
# 2. Inventory sheet - update column E (Amount)
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("Inventory")

$invChanges = @(
    @(4, "500.00000000000034"),
    @(5, "0.0"),
    @(6, "399.99999999999994"),
    @(8, "300.0000000000004"),
    @(10, "280.00000000000034"),
    @(12, "270.00000000000034"),
    @(13, "0.0"),
    @(14, "260.00000000000034"),
    @(16, "250.00000000000034"),
    @(17, "0.0"),
    @(18, "150.00000000000034"),
    @(19, "0.0"),
    @(20, "0.0"),
    @(23, "0.0"),
    @(24, "0.0"),
    @(30, "479.9999999999994"),
    @(32, "469.9999999999994"),
    @(33, "0.0"),
    @(34, "369.99999999999966"),
    @(36, "269.99999999999943"),
    @(37, "0.0"),
    @(38, "169.99999999999943"),
    @(39, "0.0"),
    @(40, "69.99999999999945"),
    @(41, "0.0"),
    @(42, "59.99999999999909"),
    @(44, "49.999999999999446"),
    @(47, "0.0"),
    @(48, "241.9999999999991"),
    @(50, "999.9999999999998"),
    @(51, "219.99999999999915"),
    @(52, "1000.0000000000005"),
    @(53, "1199.9999999999986")
)

foreach ($row in $invChanges) {
    Set-Num $wsInv $row[0] 5 $row[1]
}

# ---------------------------------------------------------------------
# 3. New sheet "RawMaterial" : Period | RawMaterial | Amount
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRM = $wb.Worksheets.Add($null, $last)
$wsRM.Name = "RawMaterial"

$wsRM.Cells.Item(1, 1).Value = "Period"
$wsRM.Cells.Item(1, 2).Value = "RawMaterial"
$wsRM.Cells.Item(1, 3).Value = "Amount"

$rmData = @(
    @(202201, "R1", "499.9999999999994"),
    @(202202, "R1", "0.0"),
    @(202203, "R1", "-2.335909243811332e-13"),
    @(202204, "R1", "0.0"),
    @(202205, "R1", "0.0"),
    @(202206, "R1", "0.0"),
    @(202207, "R1", "0.0"),
    @(202208, "R1", "2.8421709430404007e-13"),
    @(202209, "R1", "6.821210263296962e-13"),
    @(202210, "R1", "210.9999999999999"),
    @(202211, "R1", "499.9999999999999"),
    @(202212, "R1", "500.0000000000001")
)

$r = 2
foreach ($row in $rmData) {
    $wsRM.Cells.Item($r, 1).Value = $row[0]
    $wsRM.Cells.Item($r, 2).Value = $row[1]
    Set-Num $wsRM $r 3 $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 4. New sheet "RawMaterialInventory" : Period | RawMaterial | Mill | Amount
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRMI = $wb.Worksheets.Add($null, $last)
$wsRMI.Name = "RawMaterialInventory"

$wsRMI.Cells.Item(1, 1).Value = "Period"
$wsRMI.Cells.Item(1, 2).Value = "RawMaterial"
$wsRMI.Cells.Item(1, 3).Value = "Mill"
$wsRMI.Cells.Item(1, 4).Value = "Amount"

$rmiData = @(
    @(202200, "R1", "M1", "0.0"),
    @(202201, "R1", "M1", "0.0"),
    @(202202, "R1", "M1", "5.684341886080801e-13"),
    @(202203, "R1", "M1", "0.0"),
    @(202204, "R1", "M1", "0.0"),
    @(202205, "R1", "M1", "0.0"),
    @(202206, "R1", "M1", "0.0"),
    @(202207, "R1", "M1", "0.0"),
    @(202208, "R1", "M1", "0.0"),
    @(202209, "R1", "M1", "0.0"),
    @(202210, "R1", "M1", "0.0"),
    @(202211, "R1", "M1", "0.0"),
    @(202212, "R1", "M1", "0.0"),
    @(202200, "R1", "M2", "0.0"),
    @(202201, "R1", "M2", "0.0"),
    @(202202, "R1", "M2", "0.0"),
    @(202203, "R1", "M2", "0.0"),
    @(202204, "R1", "M2", "0.0"),
    @(202205, "R1", "M2", "0.0"),
    @(202206, "R1", "M2", "0.0"),
    @(202207, "R1", "M2", "0.0"),
    @(202208, "R1", "M2", "0.0"),
    @(202209, "R1", "M2", "0.0"),
    @(202210, "R1", "M2", "0.0"),
    @(202211, "R1", "M2", "0.0"),
    @(202212, "R1", "M2", "0.0")
)

$r = 2
foreach ($row in $rmiData) {
    $wsRMI.Cells.Item($r, 1).Value = $row[0]
    $wsRMI.Cells.Item($r, 2).Value = $row[1]
    $wsRMI.Cells.Item($r, 3).Value = $row[2]
    Set-Num $wsRMI $r 4 $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 5. New sheet "PMRunning" : PM | Running
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPMR = $wb.Worksheets.Add($null, $last)
$wsPMR.Name = "PMRunning"

$wsPMR.Cells.Item(1, 1).Value = "PM"
$wsPMR.Cells.Item(1, 2).Value = "Running"

$pmrData = @(
    @("PM1", "1.0"),
    @("PM2", "0.0")
)

$r = 2
foreach ($row in $pmrData) {
    $wsPMR.Cells.Item($r, 1).Value = $row[0]
    Set-Num $wsPMR $r 2 $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Restore original active sheet
# ---------------------------------------------------------------------
$wsProd.Activate()
